# Remove two "note" rows from the patient-card summary (survey sheet):
#   - row 46: note "nick"    -> "Nickname: **${aka_ctx}**"
#   - row 48: note "gender_n"-> "Gender Identity: **${gender_ctx}**"
# Deleting the entire rows shifts everything below them up, which also
# naturally updates the dependent data validation range and compacts the
# shared-string table — matching the target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Delete from the bottom up so row numbers of the two target rows don't
# shift out from under us while we work.
$ws.Rows.Item(48).Delete()
$ws.Rows.Item(46).Delete()
